# Helper: force a cell to hold literal TEXT (not an auto-converted number),
# matching the workbook's existing convention for these "numeric-looking"
# text columns, then strip the temporary "@" number-format style back off
# so no stray formatting is left behind on the cell.
function Set-TextValueAt($ws, $row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet, positioned right before "总计".
#    Easiest high-fidelity way: duplicate the "2021-Q4" sheet (same
#    column layout / header labels / styles) and then overwrite its data.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)

$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# Header row (row 1) stays identical to 2021-Q4's (基金代码/基金名称/...),
# except D1 which should read "基金规模" instead of "基金金额". Plain text
# (non-numeric-looking), so a direct .Value assignment keeps it as text
# and keeps the existing bold header style (s="2") untouched.
$new.Cells.Item(1,4).Value = "基金规模"

# Row 2: 288001 / 华夏经典配置混合
$new.Cells.Item(2,1).Value = 0
Set-TextValueAt $new 2 2 "288001"
$new.Cells.Item(2,3).Value = "华夏经典配置混合"
Set-TextValueAt $new 2 4 "18.49"
Set-TextValueAt $new 2 5 "63.85"
Set-TextValueAt $new 2 6 "5.38"
Set-TextValueAt $new 2 7 "0.9948"
$new.Cells.Item(2,8).Value = 3

# Row 3: 010363 / 信达澳银匠心臻选两年持有期混合
$new.Cells.Item(3,1).Value = 1
Set-TextValueAt $new 3 2 "010363"
$new.Cells.Item(3,3).Value = "信达澳银匠心臻选两年持有期混合"
Set-TextValueAt $new 3 4 "50.40"
Set-TextValueAt $new 3 5 "92.98"
Set-TextValueAt $new 3 6 "1.97"
Set-TextValueAt $new 3 7 "0.9929"
$new.Cells.Item(3,8).Value = 6

# Row 4: 011346 / 淳厚鑫淳一年持有期混合型证券投资基金
$new.Cells.Item(4,1).Value = 2
Set-TextValueAt $new 4 2 "011346"
$new.Cells.Item(4,3).Value = "淳厚鑫淳一年持有期混合型证券投资基金"
Set-TextValueAt $new 4 4 "5.75"
Set-TextValueAt $new 4 5 "67.80"
Set-TextValueAt $new 4 6 "2.71"
Set-TextValueAt $new 4 7 "0.1558"
$new.Cells.Item(4,8).Value = 5

# Row 5: 012454 / 淳厚鑫悦混合A
$new.Cells.Item(5,1).Value = 3
Set-TextValueAt $new 5 2 "012454"
$new.Cells.Item(5,3).Value = "淳厚鑫悦混合A"
Set-TextValueAt $new 5 4 "3.31"
Set-TextValueAt $new 5 5 "76.84"
Set-TextValueAt $new 5 6 "2.73"
Set-TextValueAt $new 5 7 "0.0904"
$new.Cells.Item(5,8).Value = 9

# Row 6 (new row, sheet only had 4 data rows before): 012455 / 淳厚鑫悦混合C
# Copy A5's style down to A6 first so the row-number column keeps the
# bold/centered look used throughout the rest of the column.
$new.Cells.Item(5,1).Copy()
$new.Cells.Item(6,1).PasteSpecial(-4122)
$new.Cells.Item(6,1).Value = 4
Set-TextValueAt $new 6 2 "012455"
$new.Cells.Item(6,3).Value = "淳厚鑫悦混合C"
Set-TextValueAt $new 6 4 "0.79"
Set-TextValueAt $new 6 5 "76.84"
Set-TextValueAt $new 6 6 "2.73"
Set-TextValueAt $new 6 7 "0.0216"
$new.Cells.Item(6,8).Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new data row for 2022-Q1
#    right under the header, pushing the existing rows down.
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows(2).Insert()

# The inserted row inherits the header's bold formatting; strip that off
# the B:D cells (they should be plain, like the rest of the data rows).
$totalWs.Range("B2:D2").ClearFormats()

# Give A2 the same style used by the rest of column A (copy from A3, the
# row that used to be row 2 before the insert).
$totalWs.Cells.Item(3,1).Copy()
$totalWs.Cells.Item(2,1).PasteSpecial(-4122)

$totalWs.Cells.Item(2,1).Value = 0
$totalWs.Cells.Item(2,2).Value = "2022-Q1"
$totalWs.Cells.Item(2,3).Value = 5
$totalWs.Cells.Item(2,4).Value = 2.26

# Column A is a running 0-based index, not a literal row copy - renumber
# the rows that got pushed down by the insert.
$totalWs.Cells.Item(3,1).Value = 1
$totalWs.Cells.Item(4,1).Value = 2
